$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("A4").Value = "Desenvolvimento de games"
$ws.Range("B4").Value = "Diurno"
$ws.Range("C4").Value = 0.29166666666666669
$ws.Range("D4").Value = 0.5
$ws.Range("C4:D4").NumberFormat = "h:mm:ss"
$ws.Range("E4").Value = "William Galvão"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("E4").Select()

$wb.Save()
